$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.700.54"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "2.528.06"
$ws.Range("E3").Value = "  -2.00%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'309.10"
$ws.Range("E5").Value = "  -1.88%  "

$ws.Range("D6").Value = "'100.84"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("E7").Value = "  -1.47%  "

$ws.Range("E9").Value = "  -2.77%  "

$ws.Range("D10").Value = "'35.80"
$ws.Range("E10").Value = "  -1.41%  "

$ws.Range("D12").Value = "'7.34"
$ws.Range("E12").Value = "  -3.26%  "

$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").Value = "2.916.50"
$ws.Range("E14").Value = "  -1.84%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.565.16"
$ws.Range("E15").Value = "  -1.04%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'15.39"
$ws.Range("E16").Value = "  -2.22%  "

$ws.Range("D17").Value = "'0.810"

$ws.Range("D18").Value = "42.690.07"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("E19").Value = "  -2.27%  "

$ws.Range("D20").Value = "'12.31"
$ws.Range("E20").Value = "  -2.92%  "

$ws.Range("D21").Value = "0.0₃0951"
$ws.Range("E21").Value = "  -2.06%  "

$ws.Range("D22").Value = "'69.63"
$ws.Range("E22").Value = "  +0.27%  "

$ws.Range("D23").Value = "'243.96"
$ws.Range("E23").Value = "  -2.57%  "

$ws.Range("E24").Value = "  -3.27%  "

$ws.Range("E25").Value = "  -3.60%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "'25.49"
$ws.Range("E27").Value = "  -6.20%  "

$ws.Range("E28").Value = "  -2.96%  "

$ws.Range("D29").Value = "'10.12"
$ws.Range("E29").Value = "  -2.11%  "

$ws.Range("D30").Value = "'38.68"
$ws.Range("E30").Value = "  -5.09%  "

$ws.Range("D31").Value = "'157.59"
$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("E32").Value = "  -2.07%  "

$ws.Range("E33").Value = "  +10.82%  "

$ws.Range("D34").Value = "'0.0785"
$ws.Range("E34").Value = "  -2.62%  "

$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("D36").Value = "'3.15"
$ws.Range("E36").Value = "  -8.46%  "

$ws.Range("E37").Value = "  -7.24%  "

$ws.Range("D38").Value = "'17.83"
$ws.Range("E38").Value = "  -5.64%  "

$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("E40").Value = "  -1.05%  "

$ws.Range("D41").Value = "'4.19"
$ws.Range("E41").Value = "  +3.27%  "

$ws.Range("D42").Value = "'21.81"
$ws.Range("E42").Value = "  -8.71%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("E44").Value = "  -1.60%  "

$ws.Range("D45").Value = "'3.29"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("D46").Value = "2.008.03"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D48").Value = "2.770.66"
$ws.Range("E48").Value = "  -1.83%  "

$ws.Range("E49").Value = "  -4.15%  "

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'79.12"
$ws.Range("E50").Value = "  -3.70%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'72.20"
$ws.Range("E51").Value = "  -3.92%  "

"Update complete"
